# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 14:27"

# --- Swap country names where rank order changed (adjacent swaps) ---
$ws.Range("A35").Value = "Kuwait"
$ws.Range("A36").Value = "Portugal"

$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"

$ws.Range("A210").Value = "Seychelles"
$ws.Range("A211").Value = "Montserrat"

$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("A214").Value = "Islas Virgenes Britanicas"

# --- Updated statistics per row (B=Casos totales, C=Nuevos casos, D=Casos activos,
#     E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2264220
$ws.Range("C4").Value = 569
$ws.Range("D4").Value = 931149
$ws.Range("E4").Value = 1212380
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 120691

# Row 7: India
$ws.Range("B7").Value = 382281
$ws.Range("C7").Value = 1190
$ws.Range("D7").Value = 205744
$ws.Range("E7").Value = 163922
$ws.Range("G7").Value = 11
$ws.Range("H7").Value = 12615

# Row 14: Alemania
$ws.Range("B14").Value = 190264
$ws.Range("C14").Value = 138
$ws.Range("E14").Value = 7212
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = 8952

# Row 27: Bielorrusia
$ws.Range("B27").Value = 57333
$ws.Range("C27").Value = 676
$ws.Range("D27").Value = 35275
$ws.Range("E27").Value = 21721
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 337

# Row 30: Paises Bajos
$ws.Range("B30").Value = 49426
$ws.Range("C30").Value = 107
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 6081

# Row 35: now Kuwait (updated data)
$ws.Range("B35").Value = 38678
$ws.Range("C35").Value = 604
$ws.Range("D35").Value = 30190
$ws.Range("E35").Value = 8175
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 313

# Row 36: now Portugal (previous Portugal data, unchanged)
$ws.Range("B36").Value = 38089
$ws.Range("D36").Value = 24010
$ws.Range("E36").Value = 12555
$ws.Range("H36").Value = 1524

# Row 50: Barein
$ws.Range("E50").Value = 5677
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 57

# Row 60: Dinamarca
$ws.Range("B60").Value = 12391
$ws.Range("C60").Value = 47
$ws.Range("D60").Value = 11282
$ws.Range("E60").Value = 509

# Row 74: Finlandia
$ws.Range("B74").Value = 7133
$ws.Range("C74").Value = 14
$ws.Range("E74").Value = 607

# Row 101: Croacia
$ws.Range("B101").Value = 2280
$ws.Range("C101").Value = 11
$ws.Range("E101").Value = 31

# Row 104: Sri Lanka
$ws.Range("B104").Value = 1948
$ws.Range("C104").Value = 2
$ws.Range("E104").Value = 516

# Row 123: Tunez
$ws.Range("B123").Value = 1146
$ws.Range("C123").Value = 14
$ws.Range("D123").Value = 1014
$ws.Range("E123").Value = 82

# Row 156: Vietnam
$ws.Range("D156").Value = 326
$ws.Range("E156").Value = 23

# Row 210: now Seychelles (previous Seychelles data)
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# Row 211: now Montserrat (previous Montserrat data)
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Row 213: now Papua Nueva Guinea (previous Papua Nueva Guinea data)
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214: now Islas Virgenes Britanicas (previous Islas Virgenes Britanicas data)
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
